$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date header in C1, matching the styling (bold/border/center) already on B1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "13-01-2023"

# Row 2: label swaps from "total" to "Alpha Acciones"; add new FCI value column C
$ws.Range("A2").Value = "Alpha Acciones"
$ws.Range("C2").Value = 9003.809999999999

# Row 3 ("avg") gains a new FCI value in column C
$ws.Range("C3").Value = 9003.809999999999

# Row 4: label swaps from "Alpha Acciones" to "total"; add new FCI value column C
$ws.Range("A4").Value = "total"
$ws.Range("C4").Value = 9003.809999999999
